# Auto-generated edit script for Chocobo_Profits workbook updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 902.9798
$ws.Range("I15").Value = 902.9798
$ws.Range("K15").Value = 2708.9394
$ws.Range("M15").Value = -2539.9394

$ws.Range("H51").Value = 15000
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()

$ws.Range("H64").Value = 3179.8
$ws.Range("I64").Value = 3149.5
$ws.Range("K64").Value = 3149.5
$ws.Range("M64").Value = -2901.5

$ws.Range("H67").Value = 3179.8
$ws.Range("I67").Value = 3149.5
$ws.Range("K67").Value = 3149.5
$ws.Range("M67").Value = -2291.5

$ws.Range("H70").Value = 4104.231
$ws.Range("I70").Value = 1975
$ws.Range("J70").Value = 5050.5557
$ws.Range("K70").Value = 5925
$ws.Range("L70").Value = 15151.6671
$ws.Range("M70").Value = -5655
$ws.Range("N70").Value = -15691.6671

$ws.Range("H73").Value = 4104.231
$ws.Range("I73").Value = 1975
$ws.Range("J73").Value = 5050.5557
$ws.Range("K73").Value = 5925
$ws.Range("L73").Value = 15151.6671
$ws.Range("M73").Value = -4989
$ws.Range("N73").Value = -17023.6671

$ws.Range("H80").Value = 682.4194
$ws.Range("I80").Value = 343.41666
$ws.Range("J80").Value = 896.5263
$ws.Range("K80").Value = 1030.24998
$ws.Range("L80").Value = 2689.5789
$ws.Range("M80").Value = -32.24998000000005
$ws.Range("N80").Value = -4685.5789

$ws.Range("H83").Value = 682.4194
$ws.Range("I83").Value = 343.41666
$ws.Range("J83").Value = 896.5263
$ws.Range("K83").Value = 3090.74994
$ws.Range("L83").Value = 8068.736699999999
$ws.Range("M83").Value = 1901.25006
$ws.Range("N83").Value = -18052.7367

$ws.Range("H118").Value = 893.2308
$ws.Range("J118").Value = 1005.5714
$ws.Range("L118").Value = 3016.7142
$ws.Range("N118").Value = -6330.7142

$ws.Range("H132").Value = 274468.3
$ws.Range("I132").Value = 4077.0881
$ws.Range("J132").Value = 3338902
$ws.Range("K132").Value = 12231.2643
$ws.Range("L132").Value = 10016706
$ws.Range("M132").Value = -9701.264299999999
$ws.Range("N132").Value = -10021766

$ws.Range("H135").Value = 403.64
$ws.Range("I135").Value = 278.6842
$ws.Range("J135").Value = 799.3333
$ws.Range("K135").Value = 2508.1578
$ws.Range("L135").Value = 7193.9997
$ws.Range("M135").Value = 26.84220000000005
$ws.Range("N135").Value = -12263.9997

$ws.Range("H137").Value = 2640.2896
$ws.Range("I137").Value = 1410.0333
$ws.Range("J137").Value = 7253.75
$ws.Range("K137").Value = 4230.0999
$ws.Range("L137").Value = 21761.25
$ws.Range("M137").Value = -1680.0999
$ws.Range("N137").Value = -26861.25

$ws.Range("H138").Value = 2105.879
$ws.Range("I138").Value = 1042.4445
$ws.Range("J138").Value = 2342.1975
$ws.Range("K138").Value = 3127.3335
$ws.Range("L138").Value = 7026.592500000001
$ws.Range("M138").Value = 2012.6665
$ws.Range("N138").Value = -17306.5925

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4747.1577
$ws.Range("I32").Value = 3759.65
$ws.Range("K32").Value = 3759.65
$ws.Range("M32").Value = -3472.65

$ws.Range("H74").Value = 7301.0713
$ws.Range("I74").Value = 7176.25
$ws.Range("K74").Value = 7176.25
$ws.Range("M74").Value = -6302.25

$ws.Range("H77").Value = 7301.0713
$ws.Range("I77").Value = 7176.25
$ws.Range("K77").Value = 35881.25
$ws.Range("M77").Value = -31513.25

$ws.Range("H80").Value = 33403.555
$ws.Range("J80").Value = 33403.555
$ws.Range("L80").Value = 33403.555
$ws.Range("N80").Value = -35399.555

$ws.Range("H83").Value = 33403.555
$ws.Range("J83").Value = 33403.555
$ws.Range("L83").Value = 100210.665
$ws.Range("N83").Value = -110194.665

$ws.Range("H87").Value = 52000
$ws.Range("J87").Value = 52000
$ws.Range("L87").Value = 52000
$ws.Range("N87").Value = -54496

$ws.Range("H90").Value = 52000
$ws.Range("J90").Value = 52000
$ws.Range("L90").Value = 156000
$ws.Range("N90").Value = -168480

$ws.Range("H110").Value = 1156.7916
$ws.Range("I110").Value = 1197.1904
$ws.Range("J110").Value = 874
$ws.Range("K110").Value = 1197.1904
$ws.Range("L110").Value = 874
$ws.Range("M110").Value = 847.8096
$ws.Range("N110").Value = -4964

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 3352685.5
$ws.Range("I7").Value = 4500
$ws.Range("J7").Value = 4022322.5
$ws.Range("K7").Value = 4500
$ws.Range("L7").Value = 4022322.5
$ws.Range("M7").Value = -4387
$ws.Range("N7").Value = -4022548.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 31253976
$ws.Range("I31").Value = 1230.1428
$ws.Range("J31").Value = 55561668
$ws.Range("K31").Value = 1230.1428
$ws.Range("L31").Value = 55561668
$ws.Range("M31").Value = -935.1428000000001
$ws.Range("N31").Value = -55562258

$ws.Range("H34").Value = 31253976
$ws.Range("I34").Value = 1230.1428
$ws.Range("J34").Value = 55561668
$ws.Range("K34").Value = 1230.1428
$ws.Range("L34").Value = 55561668
$ws.Range("M34").Value = -1028.1428
$ws.Range("N34").Value = -55562072

$ws.Range("H58").Value = 1618.8046
$ws.Range("I58").Value = 1371.183
$ws.Range("J58").Value = 5679.8
$ws.Range("K58").Value = 1371.183
$ws.Range("L58").Value = 5679.8
$ws.Range("M58").Value = -1168.183
$ws.Range("N58").Value = -6085.8

$ws.Range("H86").Value = 2042.8462
$ws.Range("I86").Value = 1373.1111
$ws.Range("J86").Value = 3549.75
$ws.Range("K86").Value = 1373.1111
$ws.Range("L86").Value = 3549.75
$ws.Range("M86").Value = -250.1111000000001
$ws.Range("N86").Value = -5795.75

$ws.Range("H89").Value = 2042.8462
$ws.Range("I89").Value = 1373.1111
$ws.Range("J89").Value = 3549.75
$ws.Range("K89").Value = 6865.5555
$ws.Range("L89").Value = 17748.75
$ws.Range("M89").Value = -1249.5555
$ws.Range("N89").Value = -28980.75

$ws.Range("H107").Value = 655.6087
$ws.Range("I107").Value = 556.8946999999999
$ws.Range("K107").Value = 556.8946999999999
$ws.Range("M107").Value = 1363.1053

$ws.Range("H132").Value = 2372.7632
$ws.Range("I132").Value = 1513.5555
$ws.Range("J132").Value = 4481.727
$ws.Range("K132").Value = 4540.666499999999
$ws.Range("L132").Value = 13445.181
$ws.Range("M132").Value = -2010.666499999999
$ws.Range("N132").Value = -18505.181

$ws.Range("H136").Value = 1618.8046
$ws.Range("I136").Value = 1371.183
$ws.Range("J136").Value = 5679.8
$ws.Range("K136").Value = 4113.549
$ws.Range("L136").Value = 17039.4
$ws.Range("M136").Value = -1563.549
$ws.Range("N136").Value = -22139.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2179.9
$ws.Range("I132").Value = 686.2857
$ws.Range("K132").Value = 6176.571300000001
$ws.Range("M132").Value = -3646.571300000001

$ws.Range("H137").Value = 2278.4707
$ws.Range("J137").Value = 3584.6667
$ws.Range("L137").Value = 10754.0001
$ws.Range("N137").Value = -20954.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 22729938
$ws.Range("I80").Value = 50002080
$ws.Range("J80").Value = 3152.1667
$ws.Range("K80").Value = 50002080
$ws.Range("L80").Value = 3152.1667
$ws.Range("M80").Value = -50001082
$ws.Range("N80").Value = -5148.1667

$ws.Range("H83").Value = 22729938
$ws.Range("I83").Value = 50002080
$ws.Range("J83").Value = 3152.1667
$ws.Range("K83").Value = 250010400
$ws.Range("L83").Value = 15760.8335
$ws.Range("M83").Value = -250005408
$ws.Range("N83").Value = -25744.8335

$ws.Range("H107").Value = 5291534.5
$ws.Range("I107").Value = 308.64706
$ws.Range("J107").Value = 27779244
$ws.Range("K107").Value = 308.64706
$ws.Range("L107").Value = 27779244
$ws.Range("M107").Value = 1611.35294
$ws.Range("N107").Value = -27783084

$ws.Range("H113").Value = 2602.75
$ws.Range("I113").Value = 1803.6666
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 1803.6666
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = 366.3334
$ws.Range("N113").Value = -9340

$ws.Range("H126").Value = 3072.0605
$ws.Range("I126").Value = 2876.253
$ws.Range("J126").Value = 4491.6665
$ws.Range("K126").Value = 8628.759
$ws.Range("L126").Value = 13474.9995
$ws.Range("M126").Value = -6158.759
$ws.Range("N126").Value = -18414.9995

$ws.Range("H134").Value = 33831.35
$ws.Range("J134").Value = 33831.35
$ws.Range("L134").Value = 101494.05
$ws.Range("N134").Value = -106564.05

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1316.8422
$ws.Range("I61").Value = 1030
$ws.Range("J61").Value = 2120
$ws.Range("K61").Value = 1030
$ws.Range("L61").Value = 2120
$ws.Range("M61").Value = -828
$ws.Range("N61").Value = -2524

$ws.Range("H82").Value = 1597.375
$ws.Range("I82").Value = 656.3
$ws.Range("J82").Value = 2269.5715
$ws.Range("K82").Value = 656.3
$ws.Range("L82").Value = 2269.5715
$ws.Range("M82").Value = -295.3
$ws.Range("N82").Value = -2991.5715

$ws.Range("H85").Value = 1597.375
$ws.Range("I85").Value = 656.3
$ws.Range("J85").Value = 2269.5715
$ws.Range("K85").Value = 656.3
$ws.Range("L85").Value = 2269.5715
$ws.Range("M85").Value = 591.7
$ws.Range("N85").Value = -4765.5715

$ws.Range("H113").Value = 1316.8422
$ws.Range("I113").Value = 1030
$ws.Range("J113").Value = 2120
$ws.Range("K113").Value = 1030
$ws.Range("L113").Value = 2120
$ws.Range("M113").Value = 1140
$ws.Range("N113").Value = -6460

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 337.42856
$ws.Range("I107").Value = 278.85715
$ws.Range("J107").Value = 396
$ws.Range("K107").Value = 836.5714499999999
$ws.Range("L107").Value = 1188
$ws.Range("M107").Value = 1083.42855
$ws.Range("N107").Value = -5028

$ws.Range("H132").Value = 6668246.5
$ws.Range("I132").Value = 928.5
$ws.Range("J132").Value = 20836296
$ws.Range("K132").Value = 2785.5
$ws.Range("L132").Value = 62508888
$ws.Range("M132").Value = -255.5
$ws.Range("N132").Value = -62513948
